# Remove duplicated teacher entries from the weekly schedule grid.
# The affected cells previously contained list-like strings showing the
# same teacher assigned to multiple class slots; they are reset to "-".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "-"
$ws.Range("C4").Value = "-"
$ws.Range("C10").Value = "-"
$ws.Range("E10").Value = "-"
$ws.Range("C16").Value = "-"
$ws.Range("E16").Value = "-"
